# Update rows 25-36 on the "Artfynd" sheet:
#  - Column A (Id) gets reshuffled among the existing rows
#  - Columns Q (Ost) / R (Nord) get updated to the rounded integer
#    coordinates that correspond to the new Id for that row
#  - Columns Z (Starttid) and AB (Sluttid) are cleared out entirely
#    (no time-of-day values remain for these rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (Id, Ost, Nord)
$updates = @(
    @{ Row = 25; Id = 112183141; Ost = 572361; Nord = 6714980 },
    @{ Row = 26; Id = 112183151; Ost = 572361; Nord = 6714978 },
    @{ Row = 27; Id = 112183148; Ost = 572357; Nord = 6714903 },
    @{ Row = 28; Id = 112183150; Ost = 572358; Nord = 6714972 },
    @{ Row = 29; Id = 112183143; Ost = 572359; Nord = 6714905 },
    @{ Row = 30; Id = 112183147; Ost = 572351; Nord = 6714915 },
    @{ Row = 31; Id = 112183134; Ost = 572354; Nord = 6714968 },
    @{ Row = 32; Id = 112183146; Ost = 572346; Nord = 6714917 },
    @{ Row = 33; Id = 112183137; Ost = 572354; Nord = 6714961 },
    @{ Row = 34; Id = 112183140; Ost = 572350; Nord = 6714962 },
    @{ Row = 35; Id = 112183145; Ost = 572351; Nord = 6714907 },
    @{ Row = 36; Id = 112183149; Ost = 572345; Nord = 6714965 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("A$r").Value = $u.Id
    $ws.Range("Q$r").Value = $u.Ost
    $ws.Range("R$r").Value = $u.Nord
    $ws.Range("Z$r").ClearContents()
    $ws.Range("AB$r").ClearContents()
}
